# Applies the scraper re-run update to the North Macedonia 1-MFL 2023-2024 sheet.
# 1) Three rows (2,3,4) get their match-data columns (F:V) cyclically rotated.
# 2) Three rows (42,44,45) get their match-data columns (F:V) cyclically rotated.
# 3) Three pairs of rows (60/61, 63/65, 75/76) get their match-data columns (F:V) swapped.
# 4) Four brand-new match rows (85-88) are appended at the bottom.
#
# NOTE: this runtime has been observed to silently drop COM calls that use
# PowerShell *named* parameters (e.g. "-Row 85") inside a function body, and
# to hang when a Range address built via string interpolation is paste-special
# targeted from inside a function. To stay on the safe/well tested path, every
# helper below takes POSITIONAL parameters only, and ranges for PasteSpecial
# are built with Cells.Item(...) instead of "A85:V85" style strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchData($Row, $Home, $HomeGoals, $Away, $AwayGoals, $HomeOpenOdds, $HomeOpenTime, $HomeCloseOdds, $HomeCloseTime, $DrawOpenOdds, $DrawOpenTime, $DrawCloseOdds, $DrawCloseTime, $AwayOpenOdds, $AwayOpenTime, $AwayCloseOdds, $AwayCloseTime, $Url) {
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenTime
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseTime
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenTime
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseTime
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenTime
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseTime
    $ws.Cells.Item($Row, 22).Value = $Url
}

function Add-MatchRow($Row, $Indice, $FechaSerial, $Home, $HomeGoals, $Away, $AwayGoals, $HomeOpenOdds, $HomeOpenTime, $HomeCloseOdds, $HomeCloseTime, $DrawOpenOdds, $DrawOpenTime, $DrawCloseOdds, $DrawCloseTime, $AwayOpenOdds, $AwayOpenTime, $AwayCloseOdds, $AwayCloseTime, $Url) {
    # Clone the formatting (borders/bold index cell, date cell number format, ...)
    # of the last existing data row (84) into the new row.
    $ws.Range("A84:V84").Copy() | Out-Null
    $target = $ws.Range($ws.Cells.Item($Row, 1), $ws.Cells.Item($Row, 22))
    $target.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = "north-macedonia"
    $ws.Cells.Item($Row, 3).Value = "1-mfl"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"
    $ws.Cells.Item($Row, 5).Value = $FechaSerial

    Set-MatchData $Row $Home $HomeGoals $Away $AwayGoals $HomeOpenOdds $HomeOpenTime $HomeCloseOdds $HomeCloseTime $DrawOpenOdds $DrawOpenTime $DrawCloseOdds $DrawCloseTime $AwayOpenOdds $AwayOpenTime $AwayCloseOdds $AwayCloseTime $Url
}

# ---------------------------------------------------------------------------
# 1) Rows 2, 3, 4 - cyclic rotation of match data (new2=old3, new3=old4, new4=old2)
# ---------------------------------------------------------------------------

Set-MatchData 2 "Voska Sport" 0 "Shkendija" 1 5.03 "06/08/2023 11:43" 6.07 "06/08/2023 16:08" 3.43 "06/08/2023 11:43" 3.38 "06/08/2023 16:11" 1.61 "06/08/2023 11:43" 1.55 "06/08/2023 16:08" "https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-shkendija-tetovo/6a4vRZCN/"

Set-MatchData 3 "Tikves" 1 "Makedonija GP" 0 2.24 "05/08/2023 05:13" 2.34 "06/08/2023 16:51" 2.84 "05/08/2023 05:13" 3.1 "06/08/2023 16:12" 2.87 "05/08/2023 05:13" 2.82 "06/08/2023 16:51" "https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-makedonija-gp/xreALDzo/"

Set-MatchData 4 "Struga" 2 "KF Gostivar" 0 1.2 "06/08/2023 11:43" 1.42 "06/08/2023 16:59" 5.7 "06/08/2023 11:43" 4.1 "06/08/2023 16:59" 10.45 "06/08/2023 11:43" 6.35 "06/08/2023 16:59" "https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-kf-gostivar/0M3rQgST/"

# ---------------------------------------------------------------------------
# 2) Rows 42, 44, 45 - cyclic rotation of match data (new42=old45, new44=old42, new45=old44)
# ---------------------------------------------------------------------------

Set-MatchData 42 "Tikves" 2 "Bregalnica Stip" 1 2.17 "23/09/2023 02:13" 2.02 "24/09/2023 14:50" 2.89 "23/09/2023 02:13" 2.81 "24/09/2023 14:50" 2.93 "23/09/2023 02:13" 3.95 "24/09/2023 14:50" "https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-bregalnica-stip/bX1AFu9m/"

Set-MatchData 44 "Shkendija" 1 "Makedonija GP" 0 1.36 "23/09/2023 02:13" 1.58 "24/09/2023 14:45" 3.89 "23/09/2023 02:13" 3.46 "24/09/2023 14:50" 6.28 "23/09/2023 02:13" 5.48 "24/09/2023 14:50" "https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-makedonija-gp/ngYicbgJ/"

Set-MatchData 45 "Voska Sport" 2 "Shkupi" 3 3.35 "24/09/2023 12:13" 4.03 "24/09/2023 14:32" 3.2 "24/09/2023 12:13" 3.35 "24/09/2023 14:32" 2.02 "24/09/2023 12:13" 1.79 "24/09/2023 14:32" "https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-shkupi/21GvL1oQ/"

# ---------------------------------------------------------------------------
# 3) Rows 60/61, 63/65, 75/76 - pairwise swap of match data
# ---------------------------------------------------------------------------

Set-MatchData 60 "Sileks" 1 "Tikves" 0 1.92 "21/10/2023 01:12" 2.01 "22/10/2023 13:52" 2.99 "21/10/2023 01:12" 3.07 "22/10/2023 13:52" 3.42 "21/10/2023 01:12" 3.55 "22/10/2023 13:52" "https://www.betexplorer.com/football/north-macedonia/1-mfl/sileks-tikves/tjl5Nf7O/"

Set-MatchData 61 "Shkendija" 1 "Brera Strumica" 1 1.51 "21/10/2023 01:12" 1.6 "22/10/2023 13:58" 3.4 "21/10/2023 01:12" 3.47 "22/10/2023 13:58" 5.12 "21/10/2023 01:12" 5.26 "22/10/2023 13:58" "https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-brera-strumica/0z5HQQYo/"

Set-MatchData 63 "Makedonija GP" 0 "Struga" 3 3.87 "21/10/2023 01:12" 4.13 "22/10/2023 13:39" 2.96 "21/10/2023 01:12" 3.07 "22/10/2023 13:39" 1.81 "21/10/2023 01:12" 1.86 "22/10/2023 13:39" "https://www.betexplorer.com/football/north-macedonia/1-mfl/makedonija-gp-struga/fgh9MEMU/"

Set-MatchData 65 "Rabotnicki" 0 "Vardar" 1 1.68 "21/10/2023 01:12" 1.8 "22/10/2023 13:53" 3.13 "21/10/2023 01:12" 3.48 "22/10/2023 13:53" 4.29 "21/10/2023 01:12" 3.8 "22/10/2023 13:53" "https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-vardar/8xjcPGxC/"

Set-MatchData 75 "Bregalnica Stip" 2 "Makedonija GP" 3 2.01 "03/11/2023 01:13" 2.24 "04/11/2023 12:51" 2.86 "03/11/2023 01:13" 2.89 "04/11/2023 12:51" 3.32 "03/11/2023 01:13" 3.2 "04/11/2023 12:51" "https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-makedonija-gp/rH5dRcPd/"

Set-MatchData 76 "Rabotnicki" 1 "Sileks" 0 2.11 "03/11/2023 01:13" 2.34 "04/11/2023 12:54" 2.83 "03/11/2023 01:13" 2.77 "04/11/2023 12:54" 3.13 "03/11/2023 01:13" 3.16 "04/11/2023 12:54" "https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-sileks/UDHE2End/"

# ---------------------------------------------------------------------------
# 4) New rows 85-88 appended at the bottom
# ---------------------------------------------------------------------------

Add-MatchRow 85 84 45255.54166666666 "Bregalnica Stip" 2 "Shkendija" 2 3.7 "25/11/2023 02:12" 3.69 "25/11/2023 12:56" 3 "25/11/2023 02:12" 3.08 "25/11/2023 12:58" 1.93 "25/11/2023 02:12" 1.96 "25/11/2023 12:58" "https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-shkendija-tetovo/44yjyZhL/"

Add-MatchRow 86 85 45256.54166666666 "Tikves" 1 "Brera Strumica" 0 2.32 "26/11/2023 02:12" 2.35 "26/11/2023 12:51" 2.8 "26/11/2023 02:12" 2.74 "26/11/2023 12:51" 3.01 "26/11/2023 02:12" 3.19 "26/11/2023 12:51" "https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-brera-strumica/nTQfzg8R/"

Add-MatchRow 87 86 45256.54166666666 "Rabotnicki" 0 "KF Gostivar" 0 2.07 "26/11/2023 02:12" 2.56 "26/11/2023 12:56" 2.96 "26/11/2023 02:12" 2.62 "26/11/2023 12:57" 3.33 "26/11/2023 02:12" 3.01 "26/11/2023 12:57" "https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-kf-gostivar/6JInxFwF/"

Add-MatchRow 88 87 45256.54166666666 "Shkupi" 1 "Makedonija GP" 0 1.37 "26/11/2023 02:12" 1.4 "26/11/2023 12:57" 4.04 "26/11/2023 02:12" 3.73 "26/11/2023 12:57" 6.95 "26/11/2023 02:12" 8.06 "26/11/2023 12:57" "https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-makedonija-gp/CdJrweO8/"

Write-Host "Update applied: rows 2-4, 42-45, 60-65, 75-76 reshuffled; rows 85-88 appended."
